# Applies the changes described by the diff between the JETT FormulaTemplate
# "before" and "after" workbook revisions:
#   - Adds a "Population Different?" label + difference formula to row 6 of
#     the "Formula Test" sheet.
#   - Adds three new template sheets: "Copy Right", "ReplaceTest" and
#     "Outside Reference".

$wb = $excel.ActiveWorkbook
$wsFormulaTest = $wb.Worksheets.Item("Formula Test")

# ---------------------------------------------------------------------------
# 1. New sheet "Copy Right" (sheet4), placed after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCopyRight = $wb.Worksheets.Add($null, $lastSheet)
$wsCopyRight.Name = "Copy Right"
$wsCopyRight.Range("A1").Value = '<jt:for start="1" end="10" var="n" copyRight="true">${n}'
$wsCopyRight.Range("A2").Value = '${2*n}'
$wsCopyRight.Range("A3").Value = '$[SUM(A1+A2)]'
$wsCopyRight.Range("A4").Value = '</jt:for>'
[void]$wsCopyRight.Range("A1:A4").Select()

# ---------------------------------------------------------------------------
# 2. New sheet "ReplaceTest" (sheet5), placed after "Copy Right".
# ---------------------------------------------------------------------------
$wsReplaceTest = $wb.Worksheets.Add($null, $wsCopyRight)
$wsReplaceTest.Name = "ReplaceTest"
$wsReplaceTest.Range("A1").Value = '<jt:for start="1" end="10" var="n">${n}'
$wsReplaceTest.Range("A2").Value = '${2*n}'
$wsReplaceTest.Range("A3").Value = '$[SUM(A1+A2)]'
$wsReplaceTest.Range("A4").Value = '</jt:for>'

# ---------------------------------------------------------------------------
# 3. New sheet "Outside Reference" (sheet6), placed after "ReplaceTest".
# ---------------------------------------------------------------------------
$wsOutsideRef = $wb.Worksheets.Add($null, $wsReplaceTest)
$wsOutsideRef.Name = "Outside Reference"
$wsOutsideRef.Range("A1").Value = '${two}'
$wsOutsideRef.Range("B1").Value = '<jt:forEach items="${primes}" var="x">${x}'
$wsOutsideRef.Range("C1").Value = '$[A1 * B1]'
$wsOutsideRef.Range("D1").Value = '<jt:forEach items="${morePrimes}" var="y">${y}'
$wsOutsideRef.Range("E1").Value = '$[A1 * B1 * D1]'
$wsOutsideRef.Range("F1").Value = '</jt:forEach></jt:forEach>'

# ---------------------------------------------------------------------------
# 4. "Formula Test" sheet (sheet1): add row 6 (written last so the two new
#    shared strings land at the end of the shared-strings table, after the
#    ones used by the three new sheets above).
# ---------------------------------------------------------------------------
$wsFormulaTest.Range("A6").Value = "Population Different?"
$wsFormulaTest.Range("C6").Value = '$[B4 <> H4]'

# Leave the original first sheet active/selected, matching the template's
# original tab selection state.
[void]$wsFormulaTest.Activate()
